# A new accelerometer sample has arrived. The sheet keeps a rolling
# window of readings (x, y, z) directly below the header row, newest
# reading on top. We push the existing readings down by one row to
# make room, write the new reading into row 2, and then drop the two
# oldest readings so the window shrinks from 21 rows of data to 20
# (sheet dimension goes from A1:C22 to A1:C21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 21

# Shift rows firstDataRow..lastDataRow down by one, working bottom-up
# so earlier values aren't overwritten before they're copied. Value2
# is used (rather than Value) so plain numbers round-trip cleanly
# without picking up any stray formatting/type along the way.
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}

# Write the freshly recorded sample into the now-vacant top data row.
$ws.Cells.Item($firstDataRow, 1).Value2 = -3.092723965644837
$ws.Cells.Item($firstDataRow, 2).Value2 = 7.026303648948669
$ws.Cells.Item($firstDataRow, 3).Value2 = -1.051015242934228

# Drop the two oldest samples, now sitting past the end of the window.
$ws.Rows.Item($lastDataRow + 2).Delete()
$ws.Rows.Item($lastDataRow + 1).Delete()
